$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.161.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").Value = "'2.600.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("D5").Value = "'312.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.11%  "

$ws.Range("D6").Value = "'98.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.24%  "

$ws.Range("D7").Value = "'0.599"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").Value = "'39.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.11%  "

$ws.Range("D11").Value = "'54.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.84%  "

$ws.Range("D12").Value = "'0.0840"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").Value = "'8.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.98%  "

$ws.Range("D14").Value = "'2.996.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("E15").Value = "  +1.06%  "

$ws.Range("D16").Value = "'2.592.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").Value = "'0.917"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.44%  "

$ws.Range("D18").Value = "'14.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").Value = "'46.210.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.88%  "

$ws.Range("E20").Value = "  +0.81%  "

$ws.Range("D21").Value = "'12.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.60%  "

$ws.Range("D22").Value = "'6.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").Value = "'289.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.32%  "

$ws.Range("D24").Value = "'72.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.79%  "

$ws.Range("D25").Value = "'3.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.93%  "

$ws.Range("D26").Value = "'2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.50%  "

$ws.Range("D27").Value = "'30.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.38%  "

$ws.Range("E29").Value = "  +1.05%  "

$ws.Range("D30").Value = "'10.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.17%  "

$ws.Range("D31").Value = "'2.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.86%  "

$ws.Range("D32").Value = "'37.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.10%  "

$ws.Range("D33").Value = "'6.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.54%  "

$ws.Range("D34").Value = "'3.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.62%  "

$ws.Range("D35").Value = "'155.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.51%  "

$ws.Range("D36").Value = "'0.0839"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.33%  "

$ws.Range("E37").Value = "  -4.25%  "

$ws.Range("E38").Value = "  -5.81%  "

$ws.Range("E39").Value = "  +4.43%  "

$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("D41").Value = "'22.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.56%  "

$ws.Range("D42").Value = "'15.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "

$ws.Range("D43").Value = "'0.0331"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.10%  "

$ws.Range("D44").Value = "'3.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.41%  "

$ws.Range("D45").Value = "'3.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.58%  "

$ws.Range("D46").Value = "'2.095.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.99%  "

$ws.Range("D47").Value = "'97.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.43%  "

$ws.Range("D48").Value = "'0.998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").Value = "'9.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.46%  "

$ws.Range("D50").Value = "'108.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "

$ws.Range("D51").Value = "'0.201"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.46%  "
